$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.193.95'
$ws.Range("D3").Value = '1.911.24'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '314.61'
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").Value = '0.5076'
$ws.Range("E7").Value = '  +0.68%  '
$ws.Range("D8").Value = '0.3926'
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("D9").Value = '0.09316'
$ws.Range("E9").Value = '  -2.92%  '
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("D11").Value = '41.96'
$ws.Range("E11").Value = '  +2.91%  '
$ws.Range("D12").Value = '6.397'
$ws.Range("E12").Value = '  -1.01%  '
$ws.Range("D13").Value = '1.920.07'
$ws.Range("E13").Value = '  +2.89%  '
$ws.Range("E14").Value = '  -0.30%  '
$ws.Range("D15").Value = '7.323'
$ws.Range("E15").Value = '  -1.06%  '
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '0.00001125'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").Value = '92.51'
$ws.Range("E18").Value = '  -0.37%  '
$ws.Range("D19").Value = '0.06620'
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").Value = '18.01'
$ws.Range("E20").Value = '  +1.89%  '
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").Value = '6.225'
$ws.Range("E22").Value = '  +0.41%  '
$ws.Range("D23").Value = '28.250.01'
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("E24").Value = '  +2.27%  '
$ws.Range("D25").Value = '2.329'
$ws.Range("E25").Value = '  +1.68%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").Value = '3.376'
$ws.Range("E26").Value = '  -0.66%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '2.593'
$ws.Range("E27").Value = '  +0.70%  '
$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D28").Value = '2.140.91'
$ws.Range("E28").Value = '  +2.65%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '21.17'
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = '157.98'
$ws.Range("E30").Value = '  -0.49%  '
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = '127.25'
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '1.106'
$ws.Range("E32").Value = '  +3.74%  '
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").Value = '0.1075'
$ws.Range("E33").Value = '  +1.13%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '5.666'
$ws.Range("E34").Value = '  +0.69%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '3.609'
$ws.Range("E35").Value = '  -0.31%  '
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").Value = '9.710'
$ws.Range("E36").Value = '  +1.92%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '0.06692'
$ws.Range("E37").Value = '  -0.98%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02437'
$ws.Range("E38").Value = '  +0.89%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = '0.2210'
$ws.Range("E39").Value = '  +0.84%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '1.243'
$ws.Range("E40").Value = '  +0.61%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '1.284'
$ws.Range("E41").Value = '  +8.43%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.6522'
$ws.Range("E42").Value = '  +2.48%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '11.57'
$ws.Range("E43").Value = '  +0.51%  '
$ws.Range("B44").Value = 'InternetComputer(DFINITY)'
$ws.Range("C44").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D44").Value = '5.025'
$ws.Range("E44").Value = '  +0.41%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.6132'
$ws.Range("E46").Value = '  +2.35%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '13.41'
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").Value = '3.725'
$ws.Range("E48").Value = '  +1.71%  '
$ws.Range("B49").Value = 'WEMIXTOKEN'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").Value = '1.291'
$ws.Range("E49").Value = '  +1.22%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '2.027'
$ws.Range("E50").Value = '  +1.10%  '
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").Value = '122.61'
$ws.Range("E51").Value = '  -0.68%  '
